# This script reproduces the Katalon "Emulator data file" commit: the
# PayNowCorp sheet gets three new test-run log entries (Result + Date)
# written into the previously-blank columns A and B of rows 2-4 (the
# header row already carries "Result"/"Date" labels in A1/B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PayNowCorp")

# Row 2: Pay Now Corporate No Emulator Data
$ws.Cells.Item(2, 1).Value = "Pass"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Value = "Thu Aug 15 22:14:55 IST 2024"
$ws.Cells.Item(2, 2).Style = "Normal"

# Row 3: Pay Now Corporate Yes Emulator Data
$ws.Cells.Item(3, 1).Value = "Pass"
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 2).Value = "Thu Aug 15 22:16:21 IST 2024"
$ws.Cells.Item(3, 2).Style = "Normal"

# Row 4: Pay Now Corporate No Emulator Data Required Fields Only
$ws.Cells.Item(4, 1).Value = "Pass"
$ws.Cells.Item(4, 1).Style = "Normal"
$ws.Cells.Item(4, 2).Value = "Thu Aug 15 22:17:32 IST 2024"
$ws.Cells.Item(4, 2).Style = "Normal"
